# Build site update for LOT2041 — the "Docentes responsaveis" (faculty) rows
# collapse into the adjoining labeled rows, and the Objetivos / Programa
# resumido / Programa / Avaliacao blocks get re-pointed to different
# (already-existing) shared strings, per the authoritative XML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) The two standalone rows that only carried professor names
#    ("1112574 - Ines Conceicao Roberto" / "1097178 - Joao Batista de
#    Almeida e Silva" in columns B/C with no column-A label) are removed.
#    Deleting them shifts every row below up by two, which realigns all of
#    the remaining column-A labels (Programa resumido:, Short syllabus:,
#    Programa:, Syllabus:, Avaliacao:, Metodo:, Criterio:, Norma de
#    recuperacao:, Bibliografia:, Requisitos:) to their new row numbers
#    with no further changes needed to those labels or the row heights.
$ws.Range("A13:A14").EntireRow.Delete()

# 2) Objetivos: row (row 10) now shows the first professor's name instead
#    of the long objectives paragraph.
$ws.Range("B10").Value = "1112574 - Inês Conceição Roberto"
$ws.Range("C10").Value = "1112574 - Inês Conceição Roberto"

# 3) Programa resumido: row (row 13 after the delete) now shows the
#    activation date instead of the Portuguese short syllabus text.
$ws.Range("B13").Value = "01/01/2019"
$ws.Range("C13").Value = "01/01/2019"

# 4) Programa: row (row 15) now shows the first professor's name instead
#    of the Portuguese full syllabus text.
$ws.Range("B15").Value = "1112574 - Inês Conceição Roberto"
$ws.Range("C15").Value = "1112574 - Inês Conceição Roberto"

# 5) Metodo: row (row 18) now shows the second professor's name instead
#    of the evaluation-method text.
$ws.Range("B18").Value = "1097178 - João Batista de Almeida e Silva"
$ws.Range("C18").Value = "1097178 - João Batista de Almeida e Silva"

# 6) Criterio: row (row 19) now shows the evaluation-method text.
$ws.Range("B19").Value = "Os alunos serão avaliados formalmente por duas provas teóricas (P1 e P2). A ponderação das notas será de 50% para cada avaliação, ou seja: Média do período letivo normal = (P1+ P2)/2."
$ws.Range("C19").Value = "Os alunos serão avaliados formalmente por duas provas teóricas (P1 e P2). A ponderação das notas será de 50% para cada avaliação, ou seja: Média do período letivo normal = (P1+ P2)/2."

# 7) Norma de recuperacao: row (row 20) now shows the passing-criterion
#    text.
$ws.Range("B20").Value = "Serão aprovados os alunos que obtiverem média do período letivo normal igual ou maior que 5."
$ws.Range("C20").Value = "Serão aprovados os alunos que obtiverem média do período letivo normal igual ou maior que 5."

# 8) Bibliografia: row (row 21) now shows the recovery-norm text instead
#    of the bibliography list.
$ws.Range("B21").Value = "Aos alunos que não obtiverem média igual ou maior que 5,0 será oferecido um programa de recuperação, que será avaliado por uma prova final (PF). Neste caso, a média final do aluno será: Média Final = (Média do período letivo normal + nota prova final) / 2. Serão aprovados os alunos que obtiverem média final igual ou maior que 5,0."
$ws.Range("C21").Value = "Aos alunos que não obtiverem média igual ou maior que 5,0 será oferecido um programa de recuperação, que será avaliado por uma prova final (PF). Neste caso, a média final do aluno será: Média Final = (Média do período letivo normal + nota prova final) / 2. Serão aprovados os alunos que obtiverem média final igual ou maior que 5,0."
